{"js": "// Apply the \"Added many more features\" edits to the Lucky Fortune Cat\n// review document using the Word JavaScript API.\n//\n// Each change is a straight text replacement, so we drive everything\n// through Body.search() (matchCase, exact phrase) + Range.insertText(\n// ..., \"Replace\"), which rewrites only the matched run's text and keeps\n// paragraph/run formatting (style, bold, italic, list bullet props) intact.\n\nconst replacements = [\n  // Title (appears twice: the Heading1 at the top and the bold line near\n  // the bottom that repeats it) -- search() finds every occurrence.\n  {\n    from: \"Play Lucky Fortune Cat for Free - Review of Red Tiger's Slot Game\",\n    to: \"Play Lucky Fortune Cat Free and Enjoy Unique Graphics\",\n  },\n  // \"What we like\" bullets\n  {\n    from: \"Unique, hand-drawn graphics\",\n    to: \"Unique hand-drawn graphics\",\n  },\n  {\n    from: \"Rare 28 pay line system\",\n    to: \"Tranquil atmosphere inspired by Chinese folklore\",\n  },\n  {\n    from: \"Golden Carp wild symbol for more wins\",\n    to: \"Rare 28 pay line system for more winning chances\",\n  },\n  {\n    from: \"Simple and easy-to-use gameplay mechanics\",\n    to: \"Intuitive gameplay mechanics\",\n  },\n  // \"What we don't like\" bullets\n  {\n    from: \"Lack of additional bonus features\",\n    to: \"Limited bonus features\",\n  },\n  {\n    from: \"Limited betting range\",\n    to: \"Lack of progressive jackpot\",\n  },\n  // Meta description (italic line)\n  {\n    from: \"Play Lucky Fortune Cat for free and read our review of Red Tiger Gaming's slot game with its unique graphics, 28 pay line system, and intuitive gameplay.\",\n    to: \"Read our review of Lucky Fortune Cat, a slot game inspired by Chinese folklore. Play for free and enjoy the unique hand-drawn graphics.\",\n  },\n];\n\nconst body = context.document.body;\n\nfor (const { from, to } of replacements) {\n  const results = body.search(from, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(to, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Added many more features\" edits to the Lucky Fortune Cat\n# review document using Word COM interop (Find/Replace).\n#\n# Each change is a straight text replacement. We drive Find.Execute with\n# MatchCase = $true so the \"Rare 28 pay line system\" (capitalized, bullet\n# text) edit doesn't clobber the several lower-case \"rare 28 pay line\n# system\" mentions elsewhere in the body copy, and Replace = 2 (wdReplaceAll)\n# so the title, which appears twice (Heading1 + bold recap line), gets\n# updated in both places.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\n# Title (appears twice: the Heading1 at the top and the bold line near the\n# bottom that repeats it) -- wdReplaceAll updates every occurrence.\nReplace-Text \"Play Lucky Fortune Cat for Free - Review of Red Tiger's Slot Game\" \"Play Lucky Fortune Cat Free and Enjoy Unique Graphics\"\n\n# \"What we like\" bullets\nReplace-Text \"Unique, hand-drawn graphics\" \"Unique hand-drawn graphics\"\nReplace-Text \"Rare 28 pay line system\" \"Tranquil atmosphere inspired by Chinese folklore\"\nReplace-Text \"Golden Carp wild symbol for more wins\" \"Rare 28 pay line system for more winning chances\"\nReplace-Text \"Simple and easy-to-use gameplay mechanics\" \"Intuitive gameplay mechanics\"\n\n# \"What we don't like\" bullets\nReplace-Text \"Lack of additional bonus features\" \"Limited bonus features\"\nReplace-Text \"Limited betting range\" \"Lack of progressive jackpot\"\n\n# Meta description (italic line)\nReplace-Text \"Play Lucky Fortune Cat for free and read our review of Red Tiger Gaming's slot game with its unique graphics, 28 pay line system, and intuitive gameplay.\" \"Read our review of Lucky Fortune Cat, a slot game inspired by Chinese folklore. Play for free and enjoy the unique hand-drawn graphics.\"\n"}
